# Update the interactive map worksheet:
# - Row 64 is overwritten with the data that used to live in row 66
#   (case -538, Malabia 964, etc.)
# - Rows 65 and 66 (the old -535 and the now-duplicated -538 rows) are deleted
# - The sheet dimension shrinks from A1:P66 to A1:P64

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete rows 65 and 66 entirely (they are no longer present after the update)
$ws.Rows.Item(66).Delete()
$ws.Rows.Item(65).Delete()

# Overwrite row 64 with the new values. Columns A, B, D and E hold values
# that look like numbers or a date ("-538", "7/31/2025", "15", "808609237")
# but must stay stored as literal text, like every other row in this sheet.
# Force a text format on just those cells first so COM does not silently
# coerce them to numeric/date values.
$numericLookingTextCols = @("A","B","D","E")
foreach ($col in $numericLookingTextCols) {
    $ws.Range($col + "64").NumberFormat = "@"
}

$ws.Range("A64").Value = "-538"
$ws.Range("B64").Value = "7/31/2025"
$ws.Range("C64").Value = "Malabia 964"
$ws.Range("D64").Value = "15"
$ws.Range("E64").Value = "808609237"
$ws.Range("F64").Value = "NEW"
$ws.Range("G64").Value = "Pendiente"
$ws.Range("H64").Value = "Cambiar poste mal estado por PRFV"
$ws.Range("I64").Value = 1
$ws.Range("J64").Value = "Cambio"
$ws.Range("K64").Value = "Sin equipos"
$ws.Range("L64").Value = "Poste"
$ws.Range("M64").Value = -58.433634
$ws.Range("N64").Value = -34.595018
$ws.Range("O64").Value = "Palermo"
$ws.Range("P64").Value = "Capital Sur"
